$d = $word.ActiveDocument

# The document contains five "<id>p008v_N</id>" paragraphs whose text was
# originally split across multiple runs (e.g. "<id>" / "p008v_1" / "</id>").
# Collapse each one back into a single run (keeping the Courier New /
# 7f6000 formatting of the opening "<id>" run) by replacing the whole
# tag text with itself via Find & Replace, which Word merges into one run.
$ids = @("p008v_1", "p008v_2", "p008v_3", "p008v_4", "p008v_5")

foreach ($id in $ids) {
    $text = "<id>" + $id + "</id>"

    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null
}
